$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the roster cells can be edited.
$ws.Unprotect()

# Row 10 = NOME (employee name), Row 11 = MAT/COD (employee number).
# Columns used for the roster (merged-cell aware), left to right.
$cols = @("G","H","I","J","K","L","M","O","P","Q","R","S","U","X","Y","Z","AA","AB","AC","AD")

$names = @(
    "ROBERTO F.DO NASCIMENTO",
    "ROBERTO CARLOS DINIZ",
    "SIVAMILTON AYOLPHI",
    "JOSE MARCIO DA SILVA JUNIOR",
    "MONA LISA ARRUDA",
    "LUZIA VASCONCELOS G DA SILVA",
    "JOSE VAGNER DA SILVA ALVES",
    "GERSON RONELLI F CARNEIRO",
    "CEIR FERNANDES DE SOUZA FILHO",
    "RODRIGO LEANDRO C. DANQUIMAIA",
    "SWELLEN NATASHA P. BARBOSA DE SIQUEIRA",
    "VICTOR DE LYRA",
    "MARCELO BITENCOURT",
    "RENATO DA SILVA MONCORES",
    "ARI DE OLIVEIRA SANTOS JUNIOR",
    "LUIZ CARLOS CHAVES DE OLIVEIRA",
    "VINICIUS FARJADO LIMA",
    "HEVERTTON NILDO M DO ROSARIO",
    "ARTUR LUCIO DUARTE NETO",
    "GABRIEL CAMARGO BATISTA*"
)

$mats = @(
    "1878",
    "2045",
    "2612",
    "3058",
    "3257",
    "3315",
    "3385",
    "3458",
    "3492",
    "3596",
    "3599",
    "3946",
    "3947",
    "3948",
    "3949",
    "3950",
    "3984",
    "4020",
    "4232",
    "4262"
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]

    $ws.Range("$col`10").Value = $names[$i]

    # Prefix with an apostrophe so the numeric-looking matricula stays a text value
    # (matching the rest of the MAT/COD row, which is stored as text).
    $ws.Range("$col`11").Value = "'" + $mats[$i]
}

$ws.Protect()
